$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, bordered, centered)
# by copying the format from the adjacent header cell H1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data values for rows 2 and 3 in columns I and J
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
